# CORE_holdings.xlsx update:
#   1) Bump the "as of" date in the confidential disclaimer text (A11) from
#      2021-05-06 to 2021-05-07.
#   2) Refresh the model holdings Weight/Percent-Change figures in D2:E8.
#
# The worksheet is protected (sheetProtection, no password known) with every
# cell locked by default, so a plain `.Value2 = ...` assignment throws
# "protected sheet" just like real Excel would. Instead of calling
# Unprotect()/Protect() (which would rewrite the <sheetProtection> element
# with different attributes than the original), we temporarily unlock just
# the cells we need, write the new values, then use Copy / PasteSpecial
# (paste formats only) from a still-pristine same-style neighbor cell to put
# each edited cell's style back exactly where it was. This keeps
# <sheetProtection .../> byte-for-byte identical to the source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-LockedCellValue($ws, $addr, $value, $formatSource) {
    $ws.Range($addr).Locked = $false
    $ws.Range($addr).Value2 = $value
    $ws.Range($formatSource).Copy()
    $ws.Range($addr).PasteSpecial($xlPasteFormats)
}

# --- 1) Disclaimer date text (A11) ------------------------------------------
# A12 keeps its original (unstyled) format throughout, so use it as the
# paste-formats source to restore A11 after the write.
$newDisclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."
Set-LockedCellValue $ws "A11" $newDisclaimer "A12"

# --- 2) Holdings Weight (D) / Percent Change (E) values ---------------------
# Process every cell except E8 using E8 (still pristine, style "1") as the
# format source, then fix up E8 last using the now-restored D2 as its source.
Set-LockedCellValue $ws "D2" 0.4984554498524808  "E8"
Set-LockedCellValue $ws "E2" 0.00796252927400487 "E8"
Set-LockedCellValue $ws "D3" 0.2451829573231916  "E8"
Set-LockedCellValue $ws "E3" 0.006829410055216645 "E8"
Set-LockedCellValue $ws "D4" 0.09615457137004187 "E8"
Set-LockedCellValue $ws "E4" 0.01204370499130869 "E8"
Set-LockedCellValue $ws "D5" 0.1032226814739828  "E8"
Set-LockedCellValue $ws "E5" 0.01083165044978895 "E8"
Set-LockedCellValue $ws "D6" 0.02995201197057159 "E8"
Set-LockedCellValue $ws "E6" 0.008351523203947941 "E8"
Set-LockedCellValue $ws "D7" 0.0270323280097314  "E8"
Set-LockedCellValue $ws "E7" 0.008713756940160522 "E8"
Set-LockedCellValue $ws "D8" 1                    "E8"
Set-LockedCellValue $ws "E8" 0.00840524841949164  "D2"
